# Scheduled data refresh: update cached market-board price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a handful of leves
# across the job sheets. Mirrors a periodic runner re-pulling prices.

$wb = $excel.ActiveWorkbook

# ALC - "Cutting Edge of Culinary Quality" (row 137)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7448.036
$ws.Range("I137").Value = 7775.222
$ws.Range("J137").Value = 6828.1055
$ws.Range("K137").Value = 23325.666
$ws.Range("L137").Value = 20484.3165
$ws.Range("M137").Value = -20775.666
$ws.Range("N137").Value = -25584.3165

# ARM - "246 Kinds of Cheese" (row 28)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 8000
$ws.Range("I28").Value = 8000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 8000
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -7808

# ARM - "Ingot We Trust" (row 32)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1687.2435
$ws.Range("I32").Value = 1469.1644
$ws.Range("J32").Value = 4871.2
$ws.Range("K32").Value = 1469.1644
$ws.Range("L32").Value = 4871.2
$ws.Range("M32").Value = -1182.1644
$ws.Range("N32").Value = -5445.2

# ARM - "Dealing with the Tough Stuff" (row 61)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3387.1035
$ws.Range("I61").Value = 2970.2307
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 2970.2307
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -2758.2307
$ws.Range("N61").Value = -7424

# ARM - "As the Bolt Flies" (row 74)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 51404.5
$ws.Range("I74").Value = 53731.684
$ws.Range("J74").Value = 36665.668
$ws.Range("K74").Value = 53731.684
$ws.Range("L74").Value = 36665.668
$ws.Range("M74").Value = -52857.684
$ws.Range("N74").Value = -38413.668

# ARM - "Heavy Metal Banned (L)" (row 77)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 51404.5
$ws.Range("I77").Value = 53731.684
$ws.Range("J77").Value = 36665.668
$ws.Range("K77").Value = 268658.42
$ws.Range("L77").Value = 183328.34
$ws.Range("M77").Value = -264290.42
$ws.Range("N77").Value = -192064.34

# ARM - "Home Cooking" (row 99)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 8000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -5005

# ARM - "Scheduled Maintenance" (row 110)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 11723.087
$ws.Range("I110").Value = 16761.084
$ws.Range("K110").Value = 16761.084
$ws.Range("M110").Value = -14716.084

# ARM - "Metal with Mettle" (row 136)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3387.1035
$ws.Range("I136").Value = 2970.2307
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 8910.6921
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -6360.6921
$ws.Range("N136").Value = -26100

# BSM - "High Steal" (row 94)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1310.8148
$ws.Range("I94").Value = 1253.6522
$ws.Range("J94").Value = 1639.5
$ws.Range("K94").Value = 1253.6522
$ws.Range("L94").Value = 1639.5
$ws.Range("M94").Value = -802.6522
$ws.Range("N94").Value = -2541.5

# CRP - "Wall Not Found" (row 31)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3728189.8
$ws.Range("I31").Value = 5065823.5
$ws.Range("J31").Value = 12541
$ws.Range("K31").Value = 5065823.5
$ws.Range("L31").Value = 12541
$ws.Range("M31").Value = -5065528.5
$ws.Range("N31").Value = -13131

# CRP - "Armoires of the Rich and Famous" (row 34)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3728189.8
$ws.Range("I34").Value = 5065823.5
$ws.Range("J34").Value = 12541
$ws.Range("K34").Value = 5065823.5
$ws.Range("L34").Value = 12541
$ws.Range("M34").Value = -5065621.5
$ws.Range("N34").Value = -12945

# CRP - "Spin It Like You Mean It" (row 52)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 64890
$ws.Range("J52").Value = 64890
$ws.Range("L52").Value = 64890
$ws.Range("N52").Value = -65478

# CRP - "Timber of Tenkonto" (row 122)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2103.1667
$ws.Range("I122").Value = 2423
$ws.Range("J122").Value = 1463.5
$ws.Range("K122").Value = 7269
$ws.Range("L122").Value = 4390.5
$ws.Range("M122").Value = -4819
$ws.Range("N122").Value = -9290.5

# CRP - "Wood You Be Quiet" (row 134)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 34977.39
$ws.Range("I134").Value = 34977.39
$ws.Range("K134").Value = 104932.17
$ws.Range("M134").Value = -102397.17

# CUL - "Creative Chocolate" (row 137)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6229.8184
$ws.Range("I137").Value = 6566.125
$ws.Range("J137").Value = 5333
$ws.Range("K137").Value = 19698.375
$ws.Range("L137").Value = 15999
$ws.Range("M137").Value = -14598.375
$ws.Range("N137").Value = -26199

# GSM - "Put the Metal to the Peddle" (row 102)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1766.3334
$ws.Range("I102").Value = 1745.091
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1745.091
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -123.0909999999999
$ws.Range("N102").Value = -5244

# GSM - "On Board for Lar" (row 132)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3166.4167
$ws.Range("I132").Value = 3173.913
$ws.Range("J132").Value = 2994
$ws.Range("K132").Value = 9521.739
$ws.Range("L132").Value = 8982
$ws.Range("M132").Value = -6991.739
$ws.Range("N132").Value = -14042

# LTW - "You Could Say It's a Moving Target" (row 68)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2574
$ws.Range("I68").Value = 2660
$ws.Range("J68").Value = 1800
$ws.Range("K68").Value = 2660
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = -1911
$ws.Range("N68").Value = -3298

# LTW - "They Call It Bloody Mary (L)" (row 71)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2574
$ws.Range("I71").Value = 2660
$ws.Range("J71").Value = 1800
$ws.Range("K71").Value = 13300
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -9556
$ws.Range("N71").Value = -16488

# LTW - "Hide to Go Seek" (row 93)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 859893.6
$ws.Range("I93").Value = 1239547
$ws.Range("J93").Value = 5673.5
$ws.Range("K93").Value = 1239547
$ws.Range("L93").Value = 5673.5
$ws.Range("M93").Value = -1238299
$ws.Range("N93").Value = -8169.5

# LTW - "Hell on Leather" (row 122)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 207047.84
$ws.Range("I122").Value = 235211.02
$ws.Range("J122").Value = 5211.6665
$ws.Range("K122").Value = 705633.0599999999
$ws.Range("L122").Value = 15634.9995
$ws.Range("M122").Value = -703183.0599999999
$ws.Range("N122").Value = -20534.9995

# WVR - "Heavy Armoire" (row 122)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4625.75
$ws.Range("I122").Value = 2318.4211
$ws.Range("J122").Value = 13393.6
$ws.Range("K122").Value = 6955.263300000001
$ws.Range("L122").Value = 40180.8
$ws.Range("M122").Value = -4505.263300000001
$ws.Range("N122").Value = -45080.8

# WVR - "Weaving the Envelope" (row 136)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2289.8667
$ws.Range("I136").Value = 2044.2963
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 6132.8889
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -3582.8889
$ws.Range("N136").Value = -18600
